# Applies the edits described in the commit diff to
# draft-gandhi-ippm-stamp-srpm-01.pptx
#
#  - Slide 10 title   : "... LM Message Format"  -> "... LM Test Packet Format"
#  - Slide 11 title   : same title text, same change (repeated slide)
#  - Slide 10 body    : "Stand-alone LM message, not tied to DM"
#                        -> "Stand-alone LM test packet, not tied to DM"
#  - Slide 10 body    : "Does not modify existing STAMP (which is for DM)
#                        procedure as different destination UDP port is used
#                        for direct-mode LM"
#                        -> "Does not modify existing STAMP procedure as
#                        different destination UDP port is used for
#                        direct-mode LM test packets"
#  - Slide 5 body     : "... the reply is required over the same path in
#                        reverse direction."
#                        -> "... the reply is required over the same path in
#                        the reverse direction."
#  - Slide 5           : delete the "Rectangle 5" shape/text box that reads
#                        "With this, the Session-Reflector node does not
#                        require any additional state for PM"

$p = $ppt.ActivePresentation

function Replace-InShape($shape, [string]$find, [string]$replace) {
    if ($shape.HasTextFrame) {
        $tf = $shape.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text.Contains($find)) {
                $tr.Replace($find, $replace) | Out-Null
            }
        }
    }
}

function Replace-InSlide($slide, [string]$find, [string]$replace) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        Replace-InShape $slide.Shapes.Item($i) $find $replace
    }
}

# --- Slide 10: "STAMP - Stand-alone Direct-mode LM..." slide ---
$slide10 = $p.Slides.Item(10)
Replace-InSlide $slide10 "STAMP - Stand-alone Direct-mode LM Message Format" "STAMP - Stand-alone Direct-mode LM Test Packet Format"
Replace-InSlide $slide10 "Stand-alone LM message, not tied to DM" "Stand-alone LM test packet, not tied to DM"
Replace-InSlide $slide10 "Does not modify existing STAMP (which is for DM) procedure as different destination UDP port is used for direct-mode LM" "Does not modify existing STAMP procedure as different destination UDP port is used for direct-mode LM test packets"

# --- Slide 11: repeats the same title text ---
$slide11 = $p.Slides.Item(11)
Replace-InSlide $slide11 "STAMP - Stand-alone Direct-mode LM Message Format" "STAMP - Stand-alone Direct-mode LM Test Packet Format"

# --- Slide 5: "STAMP - Session-Sender Control Code Field" slide ---
$slide5 = $p.Slides.Item(5)
Replace-InSlide $slide5 "Indicates that this test packet has been sent over a bidirectional path and the reply is required over the same path in reverse direction." "Indicates that this test packet has been sent over a bidirectional path and the reply is required over the same path in the reverse direction."

# Delete the "Rectangle 5" shape (with the Session-Reflector note) entirely.
for ($i = $slide5.Shapes.Count; $i -ge 1; $i--) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 5") {
        $shp.Delete() | Out-Null
    }
}
